$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# 1) Column C ("Förändrad") goes from 45207 to 45208 for every data row (rows 2-20).
for ($r = 2; $r -le 20; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45207) {
        $cell.Value = 45208
    }
}

# 2) Columns S:Y on rows 2-4 hold HYPERLINK() formulas pointing at a
#    "Logging_ANGE" folder; the folder name changed to "Logging_2260".
#    Walk every cell in S2:Y4 and patch any formula that still references
#    the old folder name, leaving everything else (labels, filenames)
#    untouched.
for ($r = 2; $r -le 4; $r++) {
    for ($c = 19; $c -le 25; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $f = $cell.Formula
        if ($f -and $f -like "*Logging_ANGE*") {
            $cell.Formula = $f -replace "Logging_ANGE", "Logging_2260"
        }
    }
}
